{"js": "// Replace the date line and each \"AxB=C\" multiplication-fact answer in the\n// document with the updated value described by the commit diff. Every\n// <w:t> value in this document is unique, so a plain exact-match\n// find & replace (via Body.search) for each old->new pair is safe and\n// will not disturb any other text or formatting.\nconst replacements = [\n  [\"2025-09-10 Wednesday\", \"2025-09-11 Thursday\"],\n  [\"29\u00d719=551\", \"18\u00d711=198\"],\n  [\"24\u00d743=1032\", \"13\u00d714=182\"],\n  [\"75\u00d726=1950\", \"48\u00d791=4368\"],\n  [\"56\u00d775=4200\", \"67\u00d711=737\"],\n  [\"35\u00d776=2660\", \"79\u00d767=5293\"],\n  [\"35\u00d731=1085\", \"24\u00d726=624\"],\n  [\"15\u00d781=1215\", \"47\u00d756=2632\"],\n  [\"27\u00d771=1917\", \"54\u00d755=2970\"],\n  [\"41\u00d781=3321\", \"16\u00d786=1376\"],\n  [\"24\u00d746=1104\", \"48\u00d742=2016\"],\n  [\"52\u00d713=676\", \"44\u00d777=3388\"],\n  [\"44\u00d743=1892\", \"76\u00d763=4788\"],\n  [\"96\u00d788=8448\", \"93\u00d756=5208\"],\n  [\"34\u00d730=1020\", \"86\u00d733=2838\"],\n  [\"29\u00d766=1914\", \"99\u00d781=8019\"],\n  [\"95\u00d713=1235\", \"33\u00d723=759\"],\n  [\"97\u00d748=4656\", \"77\u00d790=6930\"],\n  [\"29\u00d798=2842\", \"94\u00d731=2914\"],\n  [\"94\u00d717=1598\", \"86\u00d777=6622\"],\n  [\"53\u00d788=4664\", \"84\u00d737=3108\"],\n  [\"20\u00d745=900\", \"87\u00d763=5481\"],\n  [\"41\u00d763=2583\", \"96\u00d761=5856\"],\n  [\"68\u00d747=3196\", \"79\u00d787=6873\"],\n  [\"94\u00d796=9024\", \"76\u00d791=6916\"],\n  [\"65\u00d794=6110\", \"60\u00d753=3180\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and every two-digit-by-two-digit multiplication\n# \"A\u00d7B=C\" answer cell in the practice-sheet table. Every <w:t> run in this\n# document holds a unique string, so an exact-text Find/Replace (ReplaceAll)\n# for each old->new pair is safe and leaves formatting untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2025-09-10 Wednesday\"; New = \"2025-09-11 Thursday\" },\n    @{ Old = \"29\u00d719=551\"; New = \"18\u00d711=198\" },\n    @{ Old = \"24\u00d743=1032\"; New = \"13\u00d714=182\" },\n    @{ Old = \"75\u00d726=1950\"; New = \"48\u00d791=4368\" },\n    @{ Old = \"56\u00d775=4200\"; New = \"67\u00d711=737\" },\n    @{ Old = \"35\u00d776=2660\"; New = \"79\u00d767=5293\" },\n    @{ Old = \"35\u00d731=1085\"; New = \"24\u00d726=624\" },\n    @{ Old = \"15\u00d781=1215\"; New = \"47\u00d756=2632\" },\n    @{ Old = \"27\u00d771=1917\"; New = \"54\u00d755=2970\" },\n    @{ Old = \"41\u00d781=3321\"; New = \"16\u00d786=1376\" },\n    @{ Old = \"24\u00d746=1104\"; New = \"48\u00d742=2016\" },\n    @{ Old = \"52\u00d713=676\"; New = \"44\u00d777=3388\" },\n    @{ Old = \"44\u00d743=1892\"; New = \"76\u00d763=4788\" },\n    @{ Old = \"96\u00d788=8448\"; New = \"93\u00d756=5208\" },\n    @{ Old = \"34\u00d730=1020\"; New = \"86\u00d733=2838\" },\n    @{ Old = \"29\u00d766=1914\"; New = \"99\u00d781=8019\" },\n    @{ Old = \"95\u00d713=1235\"; New = \"33\u00d723=759\" },\n    @{ Old = \"97\u00d748=4656\"; New = \"77\u00d790=6930\" },\n    @{ Old = \"29\u00d798=2842\"; New = \"94\u00d731=2914\" },\n    @{ Old = \"94\u00d717=1598\"; New = \"86\u00d777=6622\" },\n    @{ Old = \"53\u00d788=4664\"; New = \"84\u00d737=3108\" },\n    @{ Old = \"20\u00d745=900\"; New = \"87\u00d763=5481\" },\n    @{ Old = \"41\u00d763=2583\"; New = \"96\u00d761=5856\" },\n    @{ Old = \"68\u00d747=3196\"; New = \"79\u00d787=6873\" },\n    @{ Old = \"94\u00d796=9024\"; New = \"76\u00d791=6916\" },\n    @{ Old = \"65\u00d794=6110\"; New = \"60\u00d753=3180\" }\n)\n\n# Word Find.Execute constants used below (as literal values, since this\n# interpreter does not pre-define the wd* enum constants):\n#   Wrap        : 1 = wdFindContinue (don't wrap past the point we started)\n#   Replace     : 2 = wdReplaceAll   (replace every match found)\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($item in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $null = $find.Execute(\n        $item.Old,    # FindText\n        $false,       # MatchCase\n        $false,       # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        $wdFindContinue, # Wrap\n        $false,       # Format\n        $item.New,    # ReplaceWith\n        $wdReplaceAll # Replace\n    )\n}\n"}
